# Add a new slide ("Resources :") with the Zenject playlist link, using
# the same "Title and Content" layout (layout index 2) as slide 3.
$p = $ppt.ActivePresentation

# There are 3 existing slides; append the new one as slide 4.
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# --- Title placeholder -------------------------------------------------
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Resources :"

# --- Body/content placeholder ------------------------------------------
$body = $s.Shapes.Item(2)
$body.Name = "Content Placeholder 4"

$tr = $body.TextFrame.TextRange
$tr.Text = "Zenject"
$tr.InsertAfter(" Playlist : https://youtube.com/playlist?list=PLKERDLXpXl_jNJPY2czQcfPXW4BJaGZc_") | Out-Null
